$d = $word.ActiveDocument

# 1) Merge the three runs that spell out the "Sau khi nghien cuu..." sentence
#    (the placeholder was split across runs as "...Ngh" + "i" + "...") back
#    into a single run by re-asserting the already-complete sentence; Word
#    collapses adjacent runs that carry identical formatting into one run.
$found1 = $d.Content.Find.Execute('Sau khi nghiên cứu, xem xét nội dung đơn kiến nghị, phản ánh của ${nguoiKienNghi}, nhận thấy không đủ điều kiện xử lý giải quyết.', $true, $false, $false, $false, $false, $true, 1, $false, 'Sau khi nghiên cứu, xem xét nội dung đơn kiến nghị, phản ánh của ${nguoiKienNghi}, nhận thấy không đủ điều kiện xử lý giải quyết.', 2)

# 2) Lengthen the dotted blank on the "Ly do:" line from 31 to 33 ellipsis
#    characters (plus one extra literal '.' before the "(3)").
$found2 = $d.Content.Find.Execute('Lý do: ………………………………………………………………………………….(3)', $true, $false, $false, $false, $false, $true, 1, $false, 'Lý do: ………………………………………………………………………………………..(3)', 2)

# 3) Nudge the size of the small freeform "underline" drawing that sits in its
#    own paragraph (not inside a table cell) up by ~0.05pt in each dimension,
#    matching the target geometry.
$shp = $d.Shapes.Item(1)
$shp.Width = 64.6
$shp.Height = 0.6
